$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.359.65"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.039.36"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.20"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.86"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.031.30"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.79"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.533.03"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.361.49"
$ws.Range("E17").Value = "  -1.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.036.83"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.87"
$ws.Range("E20").Value = "  -3.81%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.684"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.88"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  -5.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.74"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("E31").Value = "  +3.73%  "

$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "57.61"
$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.39"
$ws.Range("E34").Value = "  +3.99%  "

$ws.Range("E35").Value = "  -5.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "456.36"
$ws.Range("E37").Value = "  -1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.188.26"
$ws.Range("E38").Value = "  +2.62%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0786"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("E41").Value = "  +2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.04"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.87"
$ws.Range("E46").Value = "  +2.58%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.49"
$ws.Range("E47").Value = "  +3.28%  "

$ws.Range("E48").Value = "  +1.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.95"
$ws.Range("E49").Value = "  -2.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("E51").Value = "  +5.37%  "
